$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renames (lowercased / relabeled) ---
$ws.Range("A1").Value = "industry"
$ws.Range("B1").Value = "unit"
$ws.Range("C1").Value = "process"
$ws.Range("D1").Value = "carbon (kg CO2 eq)"
$ws.Range("E1").Value = "ced (MJ)"
$ws.Range("F1").Value = "climate change (kg CO2 eq)"
$ws.Range("G1").Value = "region"

# --- Data columns D/E/F re-derived for rows 2-43 ---
$rowData = @{
  2 = @("1.4152", "22.245421", "3.9459531e-05")
  3 = @("0.5778733333333333", "9.083546800000001", "1.6112642e-05")
  4 = @("0.23652368", "3.7178976", "6.5949077e-06")
  5 = @("0.0965805", "1.5181415", "2.6929206e-06")
  6 = @("0.14152", "2.2245421", "3.9459531e-06")
  7 = @("0.05778733333333334", "0.90835468", "1.6112642e-06")
  8 = @("0.07471068", "1.0351669", "2.0831321e-06")
  9 = @("0.007275803333333334", "0.13136651", "2.0286871e-07")
  10 = @("0.01948636", "0.30290106", "5.4333142e-07")
  11 = @("0.009766179333333333", "0.17633088", "2.7230699e-07")
  12 = @("0.04394780600000001", "0.79348897", "1.2253815e-06")
  13 = @("0.220942", "3.43438", "6.1604492e-06")
  14 = @("0.001039102733333333", "0.017745774", "2.8972941e-08")
  15 = @("0.07422162", "1.2675553", "2.0694958e-06")
  16 = @("0.08659189333333334", "1.4788145", "2.4144117e-06")
  17 = @("0.03042760533333333", "0.519642", "8.4840237e-07")
  18 = @("0.027709406", "0.47322064", "7.7261176e-07")
  19 = @("0.0002700122", "0.015340436", "7.528656600000001e-09")
  20 = @("0.019286586", "1.0957454", "5.3776119e-07")
  21 = @("0.02250101666666667", "1.2783696", "6.2738805e-07")
  22 = @("0.007906653333333333", "0.44920748", "2.2045847e-07")
  23 = @("0.007200325333333334", "0.40907828", "2.0076418e-07")
  24 = @("0.24252766", "3.8227162", "6.7623147e-06")
  25 = @("0.2637102266666667", "4.4007817", "7.3529409e-06")
  26 = @("0.3093901066666667", "5.163085", "8.626617099999999e-06")
  27 = @("0.1769419266666667", "2.7889552", "4.9336105e-06")
  28 = @("0.19252808", "3.212898", "5.3681938e-06")
  29 = @("0.2258772", "3.7694262", "6.2980556e-06")
  30 = @("0.05985599933333333", "1.5726809", "1.6689441e-06")
  31 = @("0.08035017333333333", "2.1111532", "2.240376e-06")
  32 = @("0.09433826000000001", "2.4786819", "2.630401e-06")
  33 = @("0.05985599933333333", "1.5726809", "1.6689441e-06")
  34 = @("0.08035017333333333", "2.1111532", "2.240376e-06")
  35 = @("0.09433826000000001", "2.4786819", "2.630401e-06")
  36 = @("0.02013722066666667", "0.31299712", "5.6147914e-07")
  37 = @("0.005486790333333333", "0.072039884", "1.5298627e-07")
  38 = @("0.003287408866666667", "0.043162676", "9.1661682e-08")
  39 = @("0.01287208866666667", "0.1690066", "3.5890799e-07")
  40 = @("0.011241504", "0.14759752", "3.1344296e-07")
  41 = @("0.0047754108", "0.062699688", "1.3315113e-07")
  42 = @("0.002882961333333333", "0.037852403", "8.038461100000001e-08")
  43 = @("0.001959142933333333", "0.025722949", "5.4626103e-08")
}

foreach ($r in $rowData.Keys) {
  $vals = $rowData[$r]
  $ws.Cells.Item([int]$r, 4).Value = [double]$vals[0]
  $ws.Cells.Item([int]$r, 5).Value = [double]$vals[1]
  $ws.Cells.Item([int]$r, 6).Value = [double]$vals[2]
}

# --- Cell comments on header row (also creates the legacyDrawing/vmlDrawing wiring) ---
$ws.Range("A1").AddComment("Data type: Categorical (text)") | Out-Null
$ws.Range("B1").AddComment("Data type: Various (e.g. kg, kWh)") | Out-Null
$ws.Range("C1").AddComment("Data type: Categorical (text)") | Out-Null
$ws.Range("D1").AddComment("Data type: Carbon footprint") | Out-Null
$ws.Range("E1").AddComment("Data type: Cumulative energy demand") | Out-Null
$ws.Range("F1").AddComment("Data type: Climate change impact") | Out-Null
$ws.Range("G1").AddComment("Data type: Categorical (text)") | Out-Null

Write-Output "edit complete"
